# Weekly roll of the "Fruta, Agrícola del Norte S.A. de Arica - Caqui" sheet.
# Each row's date/variety/quality/volume/price fields are replaced with the
# values that previously belonged to a different row (a weekly shift of the
# historical series), while the row's identity columns (market, region,
# product, unit, origin, kg/unit, etc.) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (previously held by row 4)
$ws.Range("D2").Value = 44313
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 1194

# New values for row 3 (previously held by row 5)
$ws.Range("D3").Value = 44342
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("S3").Value = 1361

# New values for row 4 (previously held by row 6)
$ws.Range("D4").Value = 44305
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("S4").Value = 1361

# New values for row 5 (previously held by row 8)
$ws.Range("D5").Value = 44301
$ws.Range("K5").Value = "Hachiya"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("S5").Value = 1139

# New values for row 6 (previously held by row 3)
$ws.Range("D6").Value = 44699
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("S6").Value = 1639

# Row 7 is unchanged.

# New values for row 8 (previously held by row 2)
$ws.Range("D8").Value = 45043
$ws.Range("K8").Value = "Fuyu"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 26000
$ws.Range("P8").Value = 25500
$ws.Range("S8").Value = 1417
